$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2021 column (P) to the table, copying the formatting from the
# preceding 2020 column (O) so the new cells render the same way.
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P4").Value = 2021

$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P5").Value = 80.900000000000006

$excel.CutCopyMode = $false

# Update the current selection to match the author's final selection.
$ws.Range("N10").Select()
